# Generate Report for Handback
# Update the recorded handoff/handback generation timestamps that were
# refreshed when the handback report was regenerated.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the d9b348f2 row
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-20 12:48:54"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the d9b348f2 row
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-20 12:48:49"
$zhcn.Range("K3").Value = "2016-08-20 12:49:11"

# "de-de" sheet: Correspond Handback DateTime for the d9b348f2 row
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K3").Value = "2016-08-20 12:49:18"
